$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge the TEL/CITY header box down into a new helper row --------------
# Insert a new row above the old "son1" data row (row 3). This pushes the
# "son1" row to row 4 and the "son2" row to row 5, and keeps the existing
# yellow/bordered formatting of row 2 intact while extending it downward.
$ws.Rows("3:3").Insert()

# The old "son1" sample-data row (now row 4) is no longer needed - remove it
# entirely so the sheet collapses back down.
$ws.Rows("4:4").Delete()

# The two left-over style-only rows below the data (now rows 5 and 6) are
# stray formatting remnants - clear them out completely.
$ws.Range("C5:D6").Clear()

# Match the new row 3's height to row 2's, since they are visually one block.
$ws.Rows("3:3").RowHeight = 27

# --- Build the merged "MERGED TEL & CITY" header box ------------------------
$ws.Range("C2:D3").Merge()
$ws.Range("C2:D3").Borders.LineStyle = 0
$ws.Range("C2:D3").Interior.Pattern = 1
$ws.Range("C2:D3").Interior.ThemeColor = 4
$ws.Range("C2:D3").Font.Bold = $true
$ws.Range("C2:D3").HorizontalAlignment = -4108
$ws.Range("C2:D3").VerticalAlignment = -4108
$ws.Range("C2").Value = "MERGED TEL & CITY"

# --- Normalize the remaining sample-data row to a single "final" value -----
$ws.Range("A4:E4").Value = "final"

# --- Add the new marker cell down at F9 -------------------------------------
$ws.Range("F9").Interior.Pattern = 1
$ws.Range("F9").Interior.ThemeColor = 4

# --- Update the current selection -------------------------------------------
[void]$ws.Range("B15:C15").Select()
